$wb = $excel.ActiveWorkbook

# New GUID-based file identifiers (old -> new) and new timestamps, per commit
# "Generate Report for Handoff": a fresh handoff run regenerated the source
# markdown guid, the xliff hashes, and the handoff/generation timestamps.

$newGuid = "9588c1f1-730c-4b65-acd8-0317123792bc"

$newMd   = "$newGuid.md"
$newPath = "e2e\$newGuid.md"

$newZhXlf = "$newGuid.61f7b87bd017242ac6480ba970142f69b096addb.zh-cn.xlf"
$newDeXlf = "$newGuid.61f7b87bd017242ac6480ba970142f69b096addb.de-de.xlf"

$newGenerateDate = "2016-09-07 09:25:30"
$newZhHandoffDate = "2016-09-07 09:25:15"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("B2").Value = $newPath
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = $newPath
}
$wsOverview.Range("G2").Value = $newGenerateDate

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
foreach ($hl in $wsZh.Hyperlinks) {
    $hl.TextToDisplay = $newMd
}
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandoffDate

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
foreach ($hl in $wsDe.Hyperlinks) {
    $hl.TextToDisplay = $newMd
}
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newGenerateDate
